$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade row (row 7) mirroring the existing BagTrade rows.
$ws.Range("A7").Value = 9407
$ws.Range("B7").Value = 9335.1200000000008
$ws.Range("C7").Value = 107.96
$ws.Range("D7").Value = 108.79
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.77
$ws.Range("G7").Value = 42609.488449074073
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"
$ws.Range("H7").Value = $true
